$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (C) column for rows 2-15 from 2023-10-05 (45204)
# to 2023-10-06 (45205).
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45205
}
